$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 45, column E (was already 1.5 - unchanged but keep consistent)
$ws.Range("E45").Value = 1.5

# New row 46
$ws.Range("B46").Value = "Produktansichten implementiert"
$ws.Range("C46").Value = "Jonas Alder"
$ws.Range("D46").Value = 41616
$ws.Range("E46").Value = 11

# New row 47
$ws.Range("B47").Value = "Produktansichten implementiert"
$ws.Range("C47").Value = "Jonas Alder"
$ws.Range("D47").Value = 41617
$ws.Range("E47").Value = 3

# New row 48
$ws.Range("B48").Value = "Produktansichten implementiert"
$ws.Range("C48").Value = "Jonas Alder"
$ws.Range("D48").Value = 41618
$ws.Range("E48").Value = 1.5

# New row 49
$ws.Range("B49").Value = "Testfälle erstellt, Fehlerbehebung"
$ws.Range("C49").Value = "Jonas Alder"
$ws.Range("D49").Value = 41619
$ws.Range("E49").Value = 1.5

$ws.Range("B50").Select()
